$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 43
$ws.Range("A43").Value = 'Enraged Muscles'
$ws.Range("B43").Value = 'You feel the muscles in your body contracting, your mind becomes eccentric with rage. Grow your strength by 3% for a maximum of 15% bonus.'
$ws.Range("C43").Value = 0.03
$ws.Range("M43").Value = 5
$ws.Range("N43").Value = 250

# Row 44
$ws.Range("A44").Value = 'Crafty Survival'
$ws.Range("B44").Value = 'Craft your own armour, craft your own weapons, the longer you survive the better. Increase Base Attack, Ac (Defence) and DUR by 4% for a maximum of 20%'
$ws.Range("E44").Value = 0.04
$ws.Range("J44").Value = 0.04
$ws.Range("K44").Value = 0.04
$ws.Range("M44").Value = 5
$ws.Range("N44").Value = 350
$ws.Range("O44").Value = 'Enraged Muscles'
$ws.Range("P44").Value = 2

# Row 45
$ws.Range("A45").Value = 'Hammer Crush'
$ws.Range("B45").Value = 'Feel the rage build as you crush your enemy. Gain 5% strength and 10% Base Attack for a max of 15% strength and 30% Base Attack.'
$ws.Range("C45").Value = 0.05
$ws.Range("J45").Value = 0.1
$ws.Range("M45").Value = 3
$ws.Range("N45").Value = 600
$ws.Range("O45").Value = 'Enraged Muscles'
$ws.Range("P45").Value = 4

# Row 46
$ws.Range("A46").Value = 'Shield of Purgatory'
$ws.Range("B46").Value = 'Conjure a magical shield from the depths of purgatory to give your self more defence! Increase the Base AC (defence bonus) by 6% for a total of 30%'
$ws.Range("K46").Value = 0.05
$ws.Range("M46").Value = 6
$ws.Range("N46").Value = 800
$ws.Range("O46").Value = 'Crafty Survival'
$ws.Range("P46").Value = 2

# Row 47
$ws.Range("A47").Value = 'Spirited Determination'
$ws.Range("B47").Value = 'Become determined to survive at all costs! Incrase Strength, Dexterity and Durability by 10% for a maximum of 50% while increasing Base Attack bonus by 6% for a total of 30%'
$ws.Range("C47").Value = 0.1
$ws.Range("D47").Value = 0.1
$ws.Range("E47").Value = 0.1
$ws.Range("J47").Value = 0.06
$ws.Range("M47").Value = 5
$ws.Range("N47").Value = 1200
$ws.Range("O47").Value = 'Hammer Crush'
$ws.Range("P47").Value = 3

# Row 48
$ws.Range("A48").Value = 'Blacksmiths Rage'
$ws.Range("B48").Value = 'Rage against the enemy, smasdhing the earth! Increase your Base attack by 12% for a max of 60% and your strength by 15% for a total of 75%'
$ws.Range("C48").Value = 0.15
$ws.Range("J48").Value = 0.12
$ws.Range("M48").Value = 5
$ws.Range("N48").Value = 1200
$ws.Range("O48").Value = 'Shield of Purgatory'
$ws.Range("P48").Value = 5

# Row 49
$ws.Range("A49").Value = 'A healthy body and mind'
$ws.Range("B49").Value = 'Increase your health and your strength to survive longer child! Increases strength by 12% for a total of 60%. Increases your base Healing by 2% for a total of 10% and increase your Base Ac (defence bonus) by 10% for a total of 50%'
$ws.Range("C49").Value = 0.12
$ws.Range("K49").Value = 0.1
$ws.Range("L49").Value = 0.02
$ws.Range("M49").Value = 5
$ws.Range("N49").Value = 1000
$ws.Range("O49").Value = 'Enraged Muscles'
$ws.Range("P49").Value = 5

# Row 50
$ws.Range("A50").Value = 'Last Stand'
$ws.Range("B50").Value = 'Take your last stand child! Incrwases Strength and Dur by 10% for a total of 60%. Increase your Base Attack by 12% for a total of 72%'
$ws.Range("C50").Value = 0.1
$ws.Range("E50").Value = 0.1
$ws.Range("J50").Value = 0.12
$ws.Range("M50").Value = 6
$ws.Range("N50").Value = 1500
$ws.Range("O50").Value = 'A healthy body and mind'
$ws.Range("P50").Value = 5
